$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Step 1: insert the new rows needed to make room for the new To-Do items ---
# Insert 1 row before (current) row 11 -> pushes "Axes labels"/"Fancy graph "/"Only store..." down by one.
$ws.Range("A11").EntireRow.Insert()

# Insert 2 rows before the (now shifted) "Excel" section header, which is currently at row 16.
$ws.Range("A16:A17").EntireRow.Insert()

# --- Step 2: introduce the brand new shared strings in the same order they were authored ---
# (Done, Meh, then the four new todo items in row order 11, 15, 20, 16)
$ws.Range("C4").Value = "Done"
$ws.Range("C3").Value = "Meh"
$ws.Range("B11").Value = "Graph y axis closer to values"
$ws.Range("B15").Value = "Filter data - remove garbage values"
$ws.Range("B20").Value = "Notes in excel"
$ws.Range("B16").Value = "Function stuff out"

# --- Step 3: (re)write every remaining cell with its final value / final position ---

# Header row
$ws.Range("A1").Value = "Experiment GUI"

# Setup section
$ws.Range("A2").Value = "Setup"
$ws.Range("B3").Value = "Check ports connected before experiment start button "
$ws.Range("B4").Value = "Don't store data before exp start button"
$ws.Range("B5").Value = "Instruction note"
$ws.Range("B6").Value = "Date input"
$ws.Range("B7").Value = "Open com ports at same time? Sync based on time"

# Experiment section
$ws.Range("A9").Value = "Experiment"
$ws.Range("B9").Value = "Graph pressure"
$ws.Range("C9").Value = "Done"
$ws.Range("B10").Value = "Formatting"
$ws.Range("C10").Value = "Done"
$ws.Range("B12").Value = "Axes labels"
$ws.Range("B13").Value = "Fancy graph "
$ws.Range("C13").Value = "Meh"
$ws.Range("B14").Value = "Only store data if exp running? Or graph all the time?"
$ws.Range("C14").Value = "Done"
$ws.Range("C15").Value = "Done"

# Excel section
$ws.Range("A18").Value = "Excel"
$ws.Range("B18").Value = "Header on file"
$ws.Range("C18").Value = "Done"
$ws.Range("B19").Value = "Store data with correct timestamp"

# Other section
$ws.Range("A22").Value = "Other"
$ws.Range("B22").Value = "Texts if over temp/pressure"
$ws.Range("B23").Value = "Hookup video/pic input"
$ws.Range("C23").Value = "Meh"

# Arduino Sensor Input section
$ws.Range("A27").Value = "Arduino Sensor Input"
$ws.Range("B28").Value = "Callibration adjustments"
$ws.Range("B29").Value = "Transfer of decimals "
$ws.Range("B30").Value = "Dewpoint"
$ws.Range("B31").Value = "Battery level"
$ws.Range("B32").Value = "Frequency of data sending?"
$ws.Range("C32").Value = "Meh"

# --- Step 4: update the selection to match the authored state ---
$ws.Range("F23").Select()
